$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Helper: write a literal text value into a cell, guaranteed to be stored
# as a shared string (t="s") even when the text looks like a number
# (e.g. "6", "33641"), by round-tripping it through a throw-away formula
# cell + Copy/PasteSpecial(values). A direct ".Value = <numeric-looking
# string>" assignment gets auto-coerced to a numeric cell, which is not
# what the source data represents here (these are id-like / code-like
# strings, not numbers).
function Set-TextValue($ws, $cellAddr, $text) {
    $helper = $ws.Range("ZZ1")
    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $helper.Clear()
}

# --- hotel_info sheet: fill in the review counters for row 2 ---
Set-TextValue $ws1 "G2" "6"
Set-TextValue $ws1 "I2" "6"

# --- review_info sheet: append three new review rows (2, 3, 4) ---

# Row 2
$ws2.Range("A2").Value = 39442
$ws2.Range("D2").Value = 1
Set-TextValue $ws2 "E2" "08/03/2018"
Set-TextValue $ws2 "F2" "https://www.tripadvisor.com/ShowUserReviews-g33641-d675950-r373865583-InTown_Suites_Denver_West-Sheridan_Colorado.html"
Set-TextValue $ws2 "G2" "33641"
Set-TextValue $ws2 "H2" "675950"
Set-TextValue $ws2 "I2" "373865583"
Set-TextValue $ws2 "J2" "05/17/2016"
Set-TextValue $ws2 "K2" "It's great"
Set-TextValue $ws2 "L2" "The rates are great, front desk very kind and knowledgeable. Christine is awesome!  The rest of the crew, is Great. If you need somewhere for a week or more, this is your place. Close to many different stores, right off the freeway. Service animal freindly, washer/dryer rates are reasonable."
$ws2.Range("M2").Value = 5
Set-TextValue $ws2 "N2" "April 2016"
Set-TextValue $ws2 "O2" " traveled with family"
$ws2.Range("V2").Value = 0
Set-TextValue $ws2 "Y2" "The rates are great, front desk very kind and knowledgeable. Christine is awesome!  The rest of the crew, is Great. If you need somewhere for a week or more, this is your place. Close to many different stores, right off the freeway. Service animal freindly, washer/dryer rates are reasonable."

# Row 3
$ws2.Range("A3").Value = 39442
$ws2.Range("D3").Value = 2
Set-TextValue $ws2 "E3" "08/03/2018"
Set-TextValue $ws2 "F3" "https://www.tripadvisor.com/ShowUserReviews-g33641-d675950-r277922315-InTown_Suites_Denver_West-Sheridan_Colorado.html"
Set-TextValue $ws2 "G3" "33641"
Set-TextValue $ws2 "H3" "675950"
Set-TextValue $ws2 "I3" "277922315"
Set-TextValue $ws2 "J3" "06/04/2015"
Set-TextValue $ws2 "K3" "Wonderful!"
Set-TextValue $ws2 "L3" "We stayed at this property instead of the Aurora one because of the significantly better reviews. Those were downright scary!   We had good luck with the one Indiana.  It was a wonderful stay here as well.  Thank you!!!!"
$ws2.Range("M3").Value = 5
Set-TextValue $ws2 "N3" "May 2015"
Set-TextValue $ws2 "O3" " traveled as a couple"
$ws2.Range("S3").Value = 5
$ws2.Range("U3").Value = 5
$ws2.Range("V3").Value = 0
Set-TextValue $ws2 "Y3" "We stayed at this property instead of the Aurora one because of the significantly better reviews. Those were downright scary!   We had good luck with the one Indiana.  It was a wonderful stay here as well.  Thank you!!!!"

# Row 4
$ws2.Range("A4").Value = 39442
$ws2.Range("D4").Value = 3
Set-TextValue $ws2 "E4" "08/03/2018"
Set-TextValue $ws2 "F4" "https://www.tripadvisor.com/ShowUserReviews-g33641-d675950-r197202771-InTown_Suites_Denver_West-Sheridan_Colorado.html"
Set-TextValue $ws2 "G4" "33641"
Set-TextValue $ws2 "H4" "675950"
Set-TextValue $ws2 "I4" "197202771"
Set-TextValue $ws2 "J4" "03/12/2014"
Set-TextValue $ws2 "K4" "Outstanding Accommodations"
Set-TextValue $ws2 "L4" "My experience during my stay has been terrific!  The staff is exceptionally courteous and friendly at all times. Any need or concern has always been promptly addressed as well. I have extended my stay four months as of now and couldn't be more satisfied with my room and accommodations. I would highly recommend InTown Suites to anyone looking for a very reasonably priced and comfortable hotel! :)"
$ws2.Range("M4").Value = 5
Set-TextValue $ws2 "N4" "November 2013"
Set-TextValue $ws2 "O4" " traveled solo"
$ws2.Range("P4").Value = 5
$ws2.Range("Q4").Value = 5
$ws2.Range("R4").Value = 5
$ws2.Range("S4").Value = 5
$ws2.Range("U4").Value = 5
$ws2.Range("V4").Value = 0
Set-TextValue $ws2 "Y4" "My experience during my stay has been terrific!  The staff is exceptionally courteous and friendly at all times. Any need or concern has always been promptly addressed as well. I have extended my stay four months as of now and couldn't be more satisfied with my room and accommodations. I would highly recommend InTown Suites to anyone looking for a very reasonably priced and comfortable hotel! :)"
